$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Статистика по годам" (Worksheets.Item(1))
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Update row 2 (year 2007) values
$ws1.Range("B2").Value = 55625
$ws1.Range("C2").Value = 62500
$ws1.Range("D2").Value = 4
$ws1.Range("E2").Value = 2

# Update row 3 (year 2008) values that survive, drop C3/E3 entirely
$ws1.Range("B3").Value = 43416
$ws1.Range("D3").Value = 3
$ws1.Range("C3").Clear()
$ws1.Range("E3").Clear()

# Remove rows 4-9 (years 2009-2014)
$ws1.Range("A4:A9").EntireRow.Delete()

# ---------------------------------------------------------------------------
# Sheet 2: "Статистика по городам" (Worksheets.Item(2))
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Update row 2 (Москва) values
$ws2.Range("B2").Value = 51050
# Write the percentage as literal text (not an auto-converted percent number):
# build it as a formula result then paste back as a plain value so the cell
# keeps its original "General" style instead of minting a percent format.
$ws2.Range("E2").Formula = '="71.43%"'
$ws2.Range("E2").Copy()
$ws2.Range("E2").PasteSpecial(-4163)

# Update row 3 (Санкт-Петербург) values
$ws2.Range("B3").Value = 48750
$ws2.Range("E3").Formula = '="28.57%"'
$ws2.Range("E3").Copy()
$ws2.Range("E3").PasteSpecial(-4163)

$excel.CutCopyMode = $false

# Remove rows 4-10
$ws2.Range("A4:A10").EntireRow.Delete()
